# Config.xlsx edit: add ES10 message-set row to "End Systems" sheet and
# leave the UI selection state the way the author left it when saving
# (End Systems row entered -> Settings glanced at -> Topology left active).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "End Systems" sheet: append row 14 (ES10), mirroring row 13 (ES8)
#    but with a not-yet-assigned VLID (0x0) - message count/VLID to be
#    filled in later per the commit message.
# ---------------------------------------------------------------------
$wsEnd = $wb.Worksheets.Item("End Systems")
$wsEnd.Activate()

$wsEnd.Range("A14").Value = "ES10"
$wsEnd.Range("B14").Value = 1
$wsEnd.Range("C14").Value = 0
$wsEnd.Range("D14").Value = 13
$wsEnd.Range("E14").Value = "0x0"
$wsEnd.Range("F14").Value = 0
$wsEnd.Range("G14").Value = "1s"
$wsEnd.Range("H14").Value = "1ms"
$wsEnd.Range("I14").Value = "1ms"
$wsEnd.Range("J14").Value = 0
$wsEnd.Range("K14").Value = 1183
$wsEnd.Range("L14").Value = 0
$wsEnd.Range("M14").Value = "10Mbps"
$wsEnd.Range("N14").Value = 15000

$wsEnd.Range("N14").Select()

# ---------------------------------------------------------------------
# 2. "Settings" sheet: selection returned to the top (A2), scrolled back
#    so the view no longer starts at column BD.
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Activate()
$wsSettings.Range("A2").Select()

# ---------------------------------------------------------------------
# 3. "Topology" sheet: left as the active tab, with A14:E21 selected.
# ---------------------------------------------------------------------
$wsTopology = $wb.Worksheets.Item("Topology")
$wsTopology.Activate()
$wsTopology.Range("A14:E21").Select()
